$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 92, shifting existing rows 92:117 down to 93:118
$ws.Rows.Item(92).Insert()

# Populate the newly inserted row 92 with the new weekly record
$ws.Cells.Item(92, 1).Value = 10
$ws.Cells.Item(92, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(92, 3).Value = "La Araucanía"
$ws.Cells.Item(92, 4).Value = 45093
$ws.Cells.Item(92, 5).Value = 9
$ws.Cells.Item(92, 6).Value = 100112010
$ws.Cells.Item(92, 7).Value = "Achicoria"
$ws.Cells.Item(92, 8).Value = "Sin especificar"
$ws.Cells.Item(92, 9).Value = "Primera"
$ws.Cells.Item(92, 10).Value = 85
$ws.Cells.Item(92, 11).Value = 10000
$ws.Cells.Item(92, 12).Value = 10000
$ws.Cells.Item(92, 13).Value = 10000
$ws.Cells.Item(92, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(92, 15).Value = "Región Metropolitana"
$ws.Cells.Item(92, 16).Value = 556
$ws.Cells.Item(92, 17).Value = 18
$ws.Cells.Item(92, 18).Value = "Hortaliza"
